$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes the existing 5 data rows down to rows 3-7),
# so a new "Bangalore" sample (BLT_001) becomes the first data row.
$ws.Rows.Item(2).Insert()

# Row 2 (new): BLT_001
$ws.Range("A2").Value = "BLT_001"
$ws.Range("B2").Value = 12.9716
$ws.Range("C2").Value = 77.5946

# Row 3 (was R002 / row2): now BLR_002 with updated coordinates
$ws.Range("A3").Value = "BLR_002"
$ws.Range("B3").Value = 12.9352
$ws.Range("C3").Value = 77.6146

# Row 4 (was R003 / row3): now MYS_001 with updated coordinates
$ws.Range("A4").Value = "MYS_001"
$ws.Range("B4").Value = 12.2958
$ws.Range("C4").Value = 76.6394

# Row 5 (was R004 / row4): now RUR_001 with updated coordinates
$ws.Range("A5").Value = "RUR_001"
$ws.Range("B5").Value = 13.3392
$ws.Range("C5").Value = 77.1135

# Row 6 (was R005 / row5): now DEL_001, reusing the old Delhi coordinates
$ws.Range("A6").Value = "DEL_001"
$ws.Range("B6").Value = 28.6139
$ws.Range("C6").Value = 77.209

# Row 7 (was row6, Kolkata): now MUM_001, reusing the old Mumbai coordinates
$ws.Range("A7").Value = "MUM_001"
$ws.Range("B7").Value = 19.076
$ws.Range("C7").Value = 72.8777

# Add four brand-new rows (8-11) below, carrying the same number format/style
# as the row above them.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C11").PasteSpecial(-4122)

# Row 8 (new): PUN_001
$ws.Range("A8").Value = "PUN_001"
$ws.Range("B8").Value = 18.5204
$ws.Range("C8").Value = 73.8567

# Row 9 (new): HYD_001
$ws.Range("A9").Value = "HYD_001"
$ws.Range("B9").Value = 17.385
$ws.Range("C9").Value = 78.4867

# Row 10 (new): VHN-001, reusing the old Chennai coordinates
$ws.Range("A10").Value = "VHN-001"
$ws.Range("B10").Value = 13.0827
$ws.Range("C10").Value = 80.2707

# Row 11 (new): AGRI_001
$ws.Range("A11").Value = "AGRI_001"
$ws.Range("B11").Value = 11.1271
$ws.Range("C11").Value = 78.6569
